# Refresh the cryptos price list (Price + Volume(1h) columns) with the
# latest scraped figures. Price values are stored as literal text (some
# look numeric, e.g. "1.00" or "0.0451") so NumberFormat is forced to
# Text ("@") before assignment wherever Excel would otherwise coerce the
# string into a number and silently drop formatting like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.514.16"
$ws.Range("E2").Value = "  +1.38%  "

$ws.Range("D3").Value = "2.978.34"
$ws.Range("E3").Value = "  +3.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.46"
$ws.Range("E5").Value = "  +4.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.42"
$ws.Range("E6").Value = "  +3.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.544"
$ws.Range("E7").Value = "  +0.93%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +2.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.25"
$ws.Range("E10").Value = "  +3.33%  "

$ws.Range("E11").Value = "  +0.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("E12").Value = "  +2.53%  "

$ws.Range("D13").Value = "3.448.39"
$ws.Range("E13").Value = "  +2.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.42"
$ws.Range("E14").Value = "  +1.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.56"
$ws.Range("E15").Value = "  +3.56%  "

$ws.Range("D16").Value = "2.976.22"
$ws.Range("E16").Value = "  +2.91%  "

$ws.Range("E17").Value = "  +6.73%  "

$ws.Range("D18").Value = "51.520.90"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.31"
$ws.Range("E19").Value = "  +3.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.42"
$ws.Range("E20").Value = "  +4.49%  "

$ws.Range("E21").Value = "  +2.22%  "

$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +3.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.27"
$ws.Range("E23").Value = "  +2.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.43"
$ws.Range("E24").Value = "  +2.51%  "

$ws.Range("E25").Value = "  +9.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.34"
$ws.Range("E26").Value = "  +21.61%  "

$ws.Range("E27").Value = "  +27.57%  "

$ws.Range("E28").Value = "  +16.51%  "

$ws.Range("E29").Value = "  +2.86%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").Value = "  +2.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.91"
$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.11"
$ws.Range("E33").Value = "  +4.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.06"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("E35").Value = "  -1.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0451"
$ws.Range("E36").Value = "  +8.70%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("E38").Value = "  +2.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.12"
$ws.Range("E39").Value = "  +2.14%  "

$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("E41").Value = "  +1.66%  "

$ws.Range("E42").Value = "  +4.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "125.66"
$ws.Range("E43").Value = "  +6.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.68"
$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("E45").Value = "  +21.27%  "

$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("E47").Value = "  +2.68%  "

$ws.Range("D48").Value = "2.035.57"
$ws.Range("E48").Value = "  +1.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.26"
$ws.Range("E49").Value = "  +4.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0332"
$ws.Range("E50").Value = "  +8.95%  "

$ws.Range("E51").Value = "  +4.37%  "
